$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new personnel names below the existing list.
$ws.Range("A7").Value = "王柏融"
$ws.Range("A8").Value = "金智媛"
$ws.Range("A9").Value = "姜諧潾"

# Drop the explicit cell style that was on the original rows so every
# cell (old + new) shares the default "Normal" style again.
$ws.Range("A1:A9").ClearFormats()

# Leave the selection on the last-entered cell, matching the new layout.
$ws.Range("A9").Select()
